# "Generate Report for Handback"
# Updates the Overview / zh-cn / de-de sheets of the localization-status
# workbook to reflect a completed handback: the status text moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and the
# per-language sheets get their "Latest Target File" / "Latest Handback
# File" / "Latest Handback DateTime" columns filled in (with a clickable
# link to the source markdown file in the "Latest Target File" column).

$wb = $excel.ActiveWorkbook

$mdFileName = "4da6c10e-7aed-4b24-84b5-4a4082d5794f.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d227dd316262d5232219f926d7ee9206f4fac785/e2e/4da6c10e-7aed-4b24-84b5-4a4082d5794f.md"
$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: both the zh-cn and de-de status columns flip to
# "handed back".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName)
$zhcn.Range("J2").Value = "4da6c10e-7aed-4b24-84b5-4a4082d5794f.79ea8e604bef695d2e4299c75717f09f9919d760.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-10-20 09:28:16"
$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName)
$dede.Range("J2").Value = "4da6c10e-7aed-4b24-84b5-4a4082d5794f.79ea8e604bef695d2e4299c75717f09f9919d760.de-de.xlf"
$dede.Range("K2").Value = "2016-10-20 09:28:34"
$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
